# Apply the "Cerrados" sheet addition + column cleanup to bip_center_standard_level.xlsx
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Abierto")

# ------------------------------------------------------------------
# 1) Clean up the stray, value-less formatting cells in columns K:S
#    (these only carried style "2" with no content) across the whole
#    data range of the "Abierto" sheet.
# ------------------------------------------------------------------
$ws1.Range("K1:S69").Clear()

# ------------------------------------------------------------------
# 2) Add the new "Cerrados" worksheet right after "Abierto" and
#    rebuild its header block (title + info rows + table header)
#    by copying the same cells (values, styles, merges) from Abierto.
# ------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "Cerrados"

# Keep the same compact row heights used at the top of "Abierto"
$newSheet.Rows(1).RowHeight = 11.25
$newSheet.Rows(2).RowHeight = 11.25

$headerCells = @(
    "B3","C3","D3",
    "B4","C4","D4",
    "B5","C5","D5",
    "B6","C6","D6",
    "B7","C7","D7",
    "B8","C8","D8",
    "B9","C9","D9",
    "B12","C12","D12","E12","F12","G12","H12","I12","J12"
)
foreach ($addr in $headerCells) {
    $ws1.Range($addr).Copy($newSheet.Range($addr))
}
$newSheet.Range("B3:D3").Merge()

# ------------------------------------------------------------------
# 3) Fix up the hidden autofilter defined name left over on "Abierto"
#    so it matches the (now trimmed) table range.
# ------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Abierto!_FilterDatabase") {
        $n.RefersTo = "=Abierto!`$B`$12:`$J`$67"
    }
}

$ws1.Select()
